$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the per-row IP addresses in column F (rows 2-6) with a single
# shared "127.0.0.1" value, consistent with the new config using localhost.
$ws.Range("F2:F6").Value = "127.0.0.1"

# Reflect the final selection left behind in the saved sheet (cell F14).
$ws.Range("F14").Select()
